$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.981.88'
$ws.Range("D3").Value = '1.894.53'
$ws.Range("E3").Value = '  +0.17%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").Value = "'244.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = "'0.3137"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("D9").Value = "'25.83"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.16%  '
$ws.Range("D10").Value = "'0.07273"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.55%  '
$ws.Range("D11").Value = "'0.08696"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +8.04%  '
$ws.Range("D12").Value = '2.088.23'
$ws.Range("E12").Value = '  +13.17%  '
$ws.Range("D13").Value = "'0.7750"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.57%  '
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").Value = "'94.60"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.55%  '
$ws.Range("D16").Value = "'6.215"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").Value = '30.125.13'
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("D18").Value = "'13.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = "'246.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.05%  '
$ws.Range("D20").Value = "'0.000007889"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.70%  '
$ws.Range("D21").Value = '2.274.55'
$ws.Range("E21").Value = '  +7.82%  '
$ws.Range("D22").Value = "'8.169"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = "'0.1633"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.71%  '
$ws.Range("D26").Value = "'9.530"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.44%  '
$ws.Range("D27").Value = "'163.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.03%  '
$ws.Range("D28").Value = "'18.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.82%  '
$ws.Range("D29").Value = "'2.054"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.36%  '
$ws.Range("E30").Value = '  +1.31%  '
$ws.Range("E31").Value = '  -0.01%  '
$ws.Range("D32").Value = "'4.527"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.18%  '
$ws.Range("D33").Value = "'4.136"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.82%  '
$ws.Range("D34").Value = "'0.05490"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.84%  '
$ws.Range("D35").Value = "'1.252"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.96%  '
$ws.Range("D36").Value = "'0.7556"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.56%  '
$ws.Range("D37").Value = "'1.001"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("D38").Value = "'2.691"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.74%  '
$ws.Range("D39").Value = "'0.01966"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.73%  '
$ws.Range("D40").Value = "'2.787"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.35%  '
$ws.Range("D41").Value = "'0.4525"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.36%  '
$ws.Range("D42").Value = "'74.31"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.74%  '
$ws.Range("D43").Value = '1.105.10'
$ws.Range("E43").Value = '  -3.11%  '
$ws.Range("D44").Value = "'6.073"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.81%  '
$ws.Range("D45").Value = "'0.8529"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("D46").Value = "'1.0000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("D47").Value = "'103.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.42%  '
$ws.Range("D48").Value = "'1.887"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("D49").Value = '2.171.99'
$ws.Range("E49").Value = '  +7.85%  '
$ws.Range("D50").Value = "'7.622"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.43%  '
$ws.Range("D51").Value = "'9.889"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.01%  '
